$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label from "Dad" to "All"
$ws.Range("A1").Value = "All"

# Update the work payment / expense / total values to include Cochez
$ws.Range("B2").Value = -6733.049999999999
$ws.Range("B3").Value = -7008.04
$ws.Range("B4").Value = -13741.09

# Remove the now-obsolete "Kenny" section and the Grand Total row (rows 5-9)
$ws.Range("A5:B9").EntireRow.Delete()
